$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.404.06'
$ws.Range("E2").Value = '  +2.61%  '
$ws.Range("D3").Value = '3.171.98'
$ws.Range("E3").Value = '  +1.51%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'533.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").Value = "'144.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.48%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +5.42%  '
$ws.Range("D9").Value = "'7.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").Value = '  +4.20%  '
$ws.Range("E11").Value = '  +3.80%  '
$ws.Range("D12").Value = '3.718.71'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = "'25.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("E15").Value = '  +3.04%  '
$ws.Range("D16").Value = '59.442.50'
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("D17").Value = '3.168.27'
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'377.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'8.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +3.65%  '
$ws.Range("D24").Value = "'70.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("E25").Value = '  +1.55%  '
$ws.Range("D26").Value = "'8.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +16.23%  '
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '0.0₃0896'
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("E30").Value = '  +3.35%  '
$ws.Range("D31").Value = "'6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").Value = "'5.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("D34").Value = "'6.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.29%  '
$ws.Range("D35").Value = "'156.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.26%  '
$ws.Range("E36").Value = '  +3.79%  '
$ws.Range("D37").Value = '2.740.26'
$ws.Range("E37").Value = '  +6.73%  '
$ws.Range("D38").Value = "'0.0708"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.33%  '
$ws.Range("D39").Value = "'25.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("D41").Value = "'4.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.44%  '
$ws.Range("E42").Value = '  +3.65%  '
$ws.Range("D43").Value = "'39.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("D44").Value = "'0.0289"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.75%  '
$ws.Range("B45").Value = 'RenzoRestakedETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D45").Value = '3.214.96'
$ws.Range("E45").Value = '  +1.55%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.62%  '
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").Value = "'0.1000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.27%  '
$ws.Range("D49").Value = "'20.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = "'0.767"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.98%  '
